# Working on auto setup: add a third "Phase3 / doing a lot" task row to the
# "Action list" sheet (gantt-tj3 export), with an extra "BLOCKER" reference
# back to the very first task ("Test tj3 A").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action list")

# New row 6: Container_Task, task_id, (start left blank), Effort, allocate, BLOCKER
$ws.Range("A6").Value = "Phase3"
$ws.Range("B6").Value = "doing a lot"
$ws.Range("D6").Value = "1d"
$ws.Range("E6").Value = "toC"
$ws.Range("F6").Value = "Test tj3 A"

# Leave the cursor on the last cell touched, like Excel does after data entry
$ws.Range("F6").Select()
